$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Projekt reszletek" heading: drop the stray leading-space run so
#    the paragraph reads straight "Projekt reszletek".
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    " Projekt részletek", $false, $false, $false, $false, $false,
    $true, 1, $false, "Projekt részletek", 2)

# ---------------------------------------------------------------------
# 2) Append the new "Temak" section at the end of the document body
#    (right before the trailing empty paragraph that precedes sectPr).
# ---------------------------------------------------------------------
$newBodyXml = '<w:p><w:pPr><w:pStyle w:val="Cmsor2"/></w:pPr><w:r><w:t>Témák</w:t></w:r></w:p>' +
    '<w:p/>' +
    '<w:p><w:r><w:rPr><w:rStyle w:val="Cmsor3Char"/></w:rPr><w:t>Csanádi Balázs Tóbiás</w:t></w:r></w:p>' +
    '<w:p><w:r><w:t>Cím: Synpatizer VST Plugin Hullámképzései és Oszcillátorai.</w:t></w:r></w:p>' +
    '<w:p><w:r><w:t xml:space="preserve">Backend – Synpatizer Hanghullámképzései, jelfeldolgozásai és oszcillátorai. Sound samplek kidolgozása. Hangeffektek kezelése. </w:t></w:r><w:r><w:br/></w:r></w:p>' +
    '<w:p/>' +
    '<w:p><w:r><w:rPr><w:rStyle w:val="Cmsor3Char"/></w:rPr><w:t>Fügedi Csaba</w:t></w:r></w:p>' +
    '<w:p><w:r><w:t>Cím: Synpatizer VST Plugin Kommunikációi és Jelfeldolgozásai.</w:t></w:r></w:p>' +
    '<w:p><w:r><w:t xml:space="preserve">Frontend és Backend – Synpatizer Inputkezelései, kommunikáció, Midi Interfacei. </w:t></w:r><w:r><w:t>Sound samplek kidolgozása.</w:t></w:r><w:r><w:t xml:space="preserve"> Design részegségeinek, kezelései.</w:t></w:r></w:p>'

$pkgXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
    $newBodyXml +
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$endRange = $d.Content
$endRange.Collapse(0)
$null = $endRange.InsertXML($pkgXml)
